$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.924.37'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.03%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.411.24'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.43%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.997'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '567.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.25'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.36%  '
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.525'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.394.22'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.92%  '
$ws.Range("E10").Value = '  -2.89%  '
$ws.Range("E11").Value = '  -0.63%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.05'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.13%  '
$ws.Range("E13").Value = '  -2.13%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.86'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.43%  '
$ws.Range("E15").Value = '  -0.41%  '
$ws.Range("E16").Value = '  -3.95%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.778.87'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.402.78'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.21%  '
$ws.Range("E19").Value = '  +8.17%  '
$ws.Range("E20").Value = '  -1.96%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '321.82'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.56%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.02'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.09'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.88%  '
$ws.Range("E24").Value = '  +0.02%  '
$ws.Range("E25").Value = '  -5.66%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '64.91'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '575.99'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.19'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -9.86%  '
$ws.Range("E29").Value = '  -1.23%  '
$ws.Range("E30").Value = '  -5.56%  '
$ws.Range("E31").Value = '  -2.64%  '
$ws.Range("E32").Value = '  -5.67%  '
$ws.Range("E33").Value = '  -4.83%  '
$ws.Range("E34").Value = '  -3.00%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '152.07'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.54%  '
$ws.Range("E37").Value = '  -2.29%  '
$ws.Range("E38").Value = '  -6.91%  '
$ws.Range("E39").Value = '  -3.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.15'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.57%  '
$ws.Range("E41").Value = '  -3.88%  '
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.08'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.78%  '
$ws.Range("E44").Value = '  -3.52%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.25'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -10.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '141.42'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.58%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₆0261'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.80%  '
$ws.Range("E48").Value = '  -4.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.584'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.82%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0498'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.25%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.19'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.09%  '
